$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F10"
$ws.Range("C2").Value = "F3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.887073
$ws.Range("H2").Value = 5.661219
$ws.Range("I2").Value = 0.04216693303329256
$ws.Range("J2").Value = 0.04216693303329256
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9880003333333333
$ws.Range("N2").Value = 2.964001
$ws.Range("O2").Value = 0.02006000579181712
$ws.Range("P2").Value = 0.02006000579181712
$ws.Range("Q2").Value = 1.864428753024333
$ws.Range("R2").Value = 16.779858777219
$ws.Range("S2").Value = 0.0008458689208710136
$ws.Range("T2").Value = 0.0008458689208710134

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F10"
$ws.Range("C3").Value = "F3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.887073
$ws.Range("H3").Value = 5.661219
$ws.Range("I3").Value = 0.04216693303329256
$ws.Range("J3").Value = 0.04216693303329256
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 46.102515
$ws.Range("N3").Value = 138.307545
$ws.Range("O3").Value = 0.9360489938269279
$ws.Range("P3").Value = 0.9360489938269277
$ws.Range("Q3").Value = 86.99881128859501
$ws.Range("R3").Value = 782.9893015973551
$ws.Range("S3").Value = 0.03947031523858095
$ws.Range("T3").Value = 0.03947031523858095

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "F10"
$ws.Range("C4").Value = "F3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.887073
$ws.Range("H4").Value = 5.661219
$ws.Range("I4").Value = 0.04216693303329256
$ws.Range("J4").Value = 0.04216693303329256
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6702083333333334
$ws.Range("N4").Value = 2.010625
$ws.Range("O4").Value = 0.01360767055921112
$ws.Range("P4").Value = 0.01360767055921112
$ws.Range("Q4").Value = 1.264732050208333
$ws.Range("R4").Value = 11.382588451875
$ws.Range("S4").Value = 0.000573793733209362
$ws.Range("T4").Value = 0.0005737937332093619

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "F10"
$ws.Range("C5").Value = "F3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.887073
$ws.Range("H5").Value = 5.661219
$ws.Range("I5").Value = 0.04216693303329256
$ws.Range("J5").Value = 0.04216693303329256
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.794312
$ws.Range("N5").Value = 2.382936
$ws.Range("O5").Value = 0.01612742706953525
$ws.Range("P5").Value = 0.01612742706953524
$ws.Range("Q5").Value = 1.498924728776
$ws.Range("R5").Value = 13.490322558984
$ws.Range("S5").Value = 0.0006800441372404024
$ws.Range("T5").Value = 0.0006800441372404022

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "F10"
$ws.Range("C6").Value = "F3"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.887073
$ws.Range("H6").Value = 5.661219
$ws.Range("I6").Value = 0.04216693303329256
$ws.Range("J6").Value = 0.04216693303329256
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.69721
$ws.Range("N6").Value = 2.09163
$ws.Range("O6").Value = 0.01415590275250867
$ws.Range("P6").Value = 0.01415590275250867
$ws.Range("Q6").Value = 1.31568616633
$ws.Range("R6").Value = 11.84117549697
$ws.Range("S6").Value = 0.000596911003390835
$ws.Range("T6").Value = 0.0005969110033908349

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "F10"
$ws.Range("C7").Value = "F3"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.86390766666667
$ws.Range("H7").Value = 38.591723
$ws.Range("I7").Value = 0.2874459722155911
$ws.Range("J7").Value = 0.287445972215591
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9880003333333333
$ws.Range("N7").Value = 2.964001
$ws.Range("O7").Value = 0.02006000579181712
$ws.Range("P7").Value = 0.02006000579181712
$ws.Range("Q7").Value = 12.70954506263589
$ws.Range("R7").Value = 114.385905563723
$ws.Range("S7").Value = 0.005766167867479262
$ws.Range("T7").Value = 0.005766167867479259

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "F10"
$ws.Range("C8").Value = "F3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.86390766666667
$ws.Range("H8").Value = 38.591723
$ws.Range("I8").Value = 0.2874459722155911
$ws.Range("J8").Value = 0.287445972215591
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 46.102515
$ws.Range("N8").Value = 138.307545
$ws.Range("O8").Value = 0.9360489938269279
$ws.Range("P8").Value = 0.9360489938269277
$ws.Range("Q8").Value = 593.0584961611152
$ws.Range("R8").Value = 5337.526465450035
$ws.Range("S8").Value = 0.2690635130720071
$ws.Range("T8").Value = 0.269063513072007

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "F10"
$ws.Range("C9").Value = "F3"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.86390766666667
$ws.Range("H9").Value = 38.591723
$ws.Range("I9").Value = 0.2874459722155911
$ws.Range("J9").Value = 0.287445972215591
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6702083333333334
$ws.Range("N9").Value = 2.010625
$ws.Range("O9").Value = 0.01360767055921112
$ws.Range("P9").Value = 0.01360767055921112
$ws.Range("Q9").Value = 8.621498117430557
$ws.Range("R9").Value = 77.593483056875
$ws.Range("S9").Value = 0.003911470093481916
$ws.Range("T9").Value = 0.003911470093481915

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "F10"
$ws.Range("C10").Value = "F3"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.86390766666667
$ws.Range("H10").Value = 38.591723
$ws.Range("I10").Value = 0.2874459722155911
$ws.Range("J10").Value = 0.287445972215591
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.794312
$ws.Range("N10").Value = 2.382936
$ws.Range("O10").Value = 0.01612742706953525
$ws.Range("P10").Value = 0.01612742706953524
$ws.Range("Q10").Value = 10.21795622652533
$ws.Range("R10").Value = 91.96160603872801
$ws.Range("S10").Value = 0.0046357639533386
$ws.Range("T10").Value = 0.004635763953338598

# Row 11
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = "F10"
$ws.Range("C11").Value = "F3"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 12.86390766666667
$ws.Range("H11").Value = 38.591723
$ws.Range("I11").Value = 0.2874459722155911
$ws.Range("J11").Value = 0.287445972215591
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.69721
$ws.Range("N11").Value = 2.09163
$ws.Range("O11").Value = 0.01415590275250867
$ws.Range("P11").Value = 0.01415590275250867
$ws.Range("Q11").Value = 8.968845064276668
$ws.Range("R11").Value = 80.71960557849
$ws.Range("S11").Value = 0.004069057229284217
$ws.Range("T11").Value = 0.004069057229284216

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "F10"
$ws.Range("C12").Value = "F3"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 30.001456
$ws.Range("H12").Value = 90.004368
$ws.Range("I12").Value = 0.6703870947511164
$ws.Range("J12").Value = 0.6703870947511162
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.9880003333333333
$ws.Range("N12").Value = 2.964001
$ws.Range("O12").Value = 0.02006000579181712
$ws.Range("P12").Value = 0.02006000579181712
$ws.Range("Q12").Value = 29.64144852848533
$ws.Range("R12").Value = 266.773036756368
$ws.Range("S12").Value = 0.01344796900346685
$ws.Range("T12").Value = 0.01344796900346685

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "F10"
$ws.Range("C13").Value = "F3"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 30.001456
$ws.Range("H13").Value = 90.004368
$ws.Range("I13").Value = 0.6703870947511164
$ws.Range("J13").Value = 0.6703870947511162
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 46.102515
$ws.Range("N13").Value = 138.307545
$ws.Range("O13").Value = 0.9360489938269279
$ws.Range("P13").Value = 0.9360489938269277
$ws.Range("Q13").Value = 1383.14257526184
$ws.Range("R13").Value = 12448.28317735656
$ws.Range("S13").Value = 0.6275151655163398
$ws.Range("T13").Value = 0.6275151655163397

# Row 14
$ws.Range("A14").Value = "M2"
$ws.Range("B14").Value = "F10"
$ws.Range("C14").Value = "F3"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 30.001456
$ws.Range("H14").Value = 90.004368
$ws.Range("I14").Value = 0.6703870947511164
$ws.Range("J14").Value = 0.6703870947511162
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.6702083333333334
$ws.Range("N14").Value = 2.010625
$ws.Range("O14").Value = 0.01360767055921112
$ws.Range("P14").Value = 0.01360767055921112
$ws.Range("Q14").Value = 20.10722582333334
$ws.Range("R14").Value = 180.96503241
$ws.Range("S14").Value = 0.009122406732519841
$ws.Range("T14").Value = 0.009122406732519838

# Row 15
$ws.Range("A15").Value = "M2"
$ws.Range("B15").Value = "F10"
$ws.Range("C15").Value = "F3"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 30.001456
$ws.Range("H15").Value = 90.004368
$ws.Range("I15").Value = 0.6703870947511164
$ws.Range("J15").Value = 0.6703870947511162
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.794312
$ws.Range("N15").Value = 2.382936
$ws.Range("O15").Value = 0.01612742706953525
$ws.Range("P15").Value = 0.01612742706953524
$ws.Range("Q15").Value = 23.830516518272
$ws.Range("R15").Value = 214.474648664448
$ws.Range("S15").Value = 0.01081161897895624
$ws.Range("T15").Value = 0.01081161897895624

# Row 16
$ws.Range("A16").Value = "M2"
$ws.Range("B16").Value = "F10"
$ws.Range("C16").Value = "F3"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 30.001456
$ws.Range("H16").Value = 90.004368
$ws.Range("I16").Value = 0.6703870947511164
$ws.Range("J16").Value = 0.6703870947511162
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.69721
$ws.Range("N16").Value = 2.09163
$ws.Range("O16").Value = 0.01415590275250867
$ws.Range("P16").Value = 0.01415590275250867
$ws.Range("Q16").Value = 20.91731513776
$ws.Range("R16").Value = 188.25583623984
$ws.Range("S16").Value = 0.009489934519833619
$ws.Range("T16").Value = 0.009489934519833617
